$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features11 = "11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), #digit/#ascii, %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, last_character_ascii, last_character_digit"
$model = "Neuron Network"
$filter0 = "0 filters: "

$rows = @(
    @{
        Row = 25
        A = "20160415_171314"
        B = 2025.436
        C = 'space after punctuation, remove multiple spaces, remove break line, trim "space" and ",", convert unicode to ascii, convert to lower'
        D = $features11
        E = $model
        F = "3 layers: [100-Sigmoid, 3-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
        G = 0.996666666666667
        H = 0.940594059405941
        I = $filter0
        J = 0.25609756097561
    },
    @{
        Row = 26
        A = "20160415_174700"
        B = 2135.522
        C = 'space after punctuation, remove multiple spaces, remove break line, trim "space" and ",", convert unicode to ascii, convert to lower'
        D = $features11
        E = $model
        F = "3 layers: [100-Sigmoid, 3-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
        G = 0.998
        H = 0.95049504950495
        I = $filter0
        J = 0.317647058823529
    },
    @{
        Row = 27
        A = "20160415_182235"
        B = 1930.165
        C = 'space after punctuation, remove multiple spaces, remove break line, trim "space" and ",", convert unicode to ascii, convert to lower'
        D = $features11
        E = $model
        F = "3 layers: [100-Sigmoid, 3-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
        G = 0.998666666666667
        H = 0.940594059405941
        I = $filter0
        J = 0.548780487804878
    }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
}
